$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "61.272.23"
$ws.Range("E2").Value = "  +0.29%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.923.93"
$ws.Range("E3").Value = "  -0.13%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.999"

# Row 5 - BNB
$ws.Range("D5").Value = "596.64"
$ws.Range("E5").Value = "  +0.45%  "

# Row 6 - Solana
$ws.Range("D6").Value = "144.88"
$ws.Range("E6").Value = "  -1.01%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.05%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -1.14%  "

# Row 9 - Toncoin
$ws.Range("E9").Value = "  +0.91%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -2.68%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -0.77%  "

# Row 12 - ShibaInu
$ws.Range("E12").Value = "  -1.24%  "

# Row 13 - Avalanche
$ws.Range("E13").Value = "  -1.26%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +0.26%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.406.55"
$ws.Range("E15").Value = "  -0.18%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "61.257.53"

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.921.28"
$ws.Range("E17").Value = "  -0.29%  "

# Row 18 - Polkadot
$ws.Range("D18").Value = "6.67"
$ws.Range("E18").Value = "  -0.75%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "430.86"
$ws.Range("E19").Value = "  -0.26%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -0.21%  "

# Row 21 - Polygon
$ws.Range("D21").Value = "0.674"
$ws.Range("E21").Value = "  -1.33%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "7.06"

# Row 23 - Litecoin
$ws.Range("D23").Value = "81.73"
$ws.Range("E23").Value = "  +0.46%  "

# Row 24 - RenderToken
$ws.Range("D24").Value = "10.85"
$ws.Range("E24").Value = "  -2.15%  "

# Row 25 - Fetch.AI
$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  -2.26%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "11.69"
$ws.Range("E26").Value = "  -2.86%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.03%  "

# Row 28 - ImmutableX
$ws.Range("E28").Value = "  -4.78%  "

# Row 29 - PancakeSwap
$ws.Range("D29").Value = "2.60"
$ws.Range("E29").Value = "  -0.73%  "

# Row 30 - NEARProtocol
$ws.Range("D30").Value = "6.90"
$ws.Range("E30").Value = "  -2.85%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +1.17%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "26.57"
$ws.Range("E32").Value = "  +0.28%  "

# Row 33 - FirstDigitalUSD
$ws.Range("E33").Value = "  +0.08%  "

# Row 34 - PEPE
$ws.Range("D34").Value = "0.0₃0880"
$ws.Range("E34").Value = "  +2.83%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  -0.38%  "

# Row 36 - Filecoin
$ws.Range("E36").Value = "  -0.64%  "

# Row 37 - dogwifhat
$ws.Range("E37").Value = "  -3.11%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -2.05%  "

# Row 40 - Cosmos
$ws.Range("D40").Value = "8.54"
$ws.Range("E40").Value = "  -0.74%  "

# Row 41 - Arweave
$ws.Range("D41").Value = "42.14"
$ws.Range("E41").Value = "  +6.30%  "

# Row 42 - TheGraph
$ws.Range("D42").Value = "0.280"
$ws.Range("E42").Value = "  -2.32%  "

# Row 43 - VeChain
$ws.Range("D43").Value = "0.0344"
$ws.Range("E43").Value = "  -0.58%  "

# Row 44 - Maker
$ws.Range("D44").Value = "2.697.29"
$ws.Range("E44").Value = "  -0.94%  "

# Row 45 - Monero
$ws.Range("D45").Value = "133.62"
$ws.Range("E45").Value = "  +2.15%  "

# Row 46 - Bittensor
$ws.Range("D46").Value = "360.07"
$ws.Range("E46").Value = "  -4.29%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "23.52"
$ws.Range("E48").Value = "  -2.72%  "

# Row 50 - ThetaToken
$ws.Range("E50").Value = "  -1.87%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  -2.47%  "
